{"js": "// 1) The \"_GoBack\" bookmark that sat at the end of the first (\"fsck\") paragraph\n//    is removed from there ...\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) ... and re-inserted a few paragraphs down, right after \"Making sure t\"\n//    (splitting \"Making sure the free block list contains ALL of the free\n//    blocks\" into two runs), marking where the author's cursor ended up.\nconst body = context.document.body;\nconst bookmarkAnchor = body.search(\"Making sure t\", { matchCase: true });\nbookmarkAnchor.load(\"items\");\nawait context.sync();\n\nif (bookmarkAnchor.items.length > 0) {\n  const afterT = bookmarkAnchor.items[0].getRange(\"End\");\n  afterT.insertBookmark(\"_GoBack\");\n}\n\n// 3) The \"The DeviceID is correct\" bullet gets struck through (paragraph\n//    mark + every run in it).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(p => p.text === \"The DeviceID is correct\");\nif (target) {\n  target.font.strikeThrough = true;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Strike through the \"The DeviceID is correct\" bullet (paragraph mark +\n#    every run in it).\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"The DeviceID is correct\") {\n        $p.Range.Font.StrikeThrough = 1\n    }\n}\n\n# 2) Move the \"_GoBack\" bookmark: it used to sit at the end of the first\n#    (\"fsck\") paragraph; it now sits right after \"Making sure t\", splitting\n#    \"Making sure the free block list contains ALL of the free blocks\" into\n#    two runs. Re-adding a bookmark with the same name relocates it.\n$target = $d.Content\n$target.Find.Execute(\"Making sure t\") | Out-Null\n$target.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $target) | Out-Null\n"}
